# sp_Blitz Check ID List - add new checks (#115 make_parallel, #116 backup
# compression default off, #117 memory pressure / forced grants, #114 NUMA
# config hardware row), per commit message:
#   "Adding #115, #182, #183, #189, #180
#    New checks for forced grants, backup compression default off, Adam
#    Machanic's make_parallel function, NUMA nodes, fixing named instance
#    check."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 179: Check #114 - Hardware - NUMA Config (no hyperlink) ---
$ws.Range("A179").Value = 114
$ws.Range("B179").Value = 250
$ws.Range("C179").Value = "Server Info"
$ws.Range("D179").Value = "Hardware - NUMA Config"

# --- Row 180: Check #115 - Parallelism Rocket Surgery (make_parallel) ---
$ws.Range("A180").Value = 115
$ws.Range("B180").Value = 110
$ws.Range("C180").Value = "Performance"
$ws.Range("D180").Value = "Parallelism Rocket Surgery"
$ws.Hyperlinks.Add($ws.Range("E180"), "http://BrentOzar.com/go/makeparallel") | Out-Null

# --- Row 181: Check #116 - Backup Compression Default Off ---
$ws.Range("A181").Value = 116
$ws.Range("B181").Value = 200
$ws.Range("C181").Value = "Informational"
$ws.Range("D181").Value = "Backup Compression Default Off"
$ws.Hyperlinks.Add($ws.Range("E181"), "http://BrentOzar.com/go/backup") | Out-Null

# --- Row 182: Check #117 - Memory Pressure Affecting Queries (forced grants) ---
$ws.Range("A182").Value = 117
$ws.Range("B182").Value = 100
$ws.Range("C182").Value = "Performance"
$ws.Range("D182").Value = "Memory Pressure Affecting Queries"
$ws.Hyperlinks.Add($ws.Range("E182"), "http://BrentOzar.com/go/grants") | Out-Null

# Leave the selection where the author ended up: one row below the newly
# added data.
$ws.Range("E183").Select() | Out-Null
